$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update nombre_aides (C) and montant_total (D) for the "Bourgogne-Franche-Comte"
# through "Provence-Alpes-Cote d'Azur" rows with the 2020-12-30 data refresh.
# Values are leading-apostrophe ('-prefixed) so Excel keeps them as text
# (matching the source workbook, where these columns are stored as strings),
# instead of being auto-coerced to numbers.
$ws.Range("C9").Value = "'316"
$ws.Range("D9").Value = "'1042987.10"
$ws.Range("C11").Value = "'532"
$ws.Range("D11").Value = "'3836269.21"
$ws.Range("C17").Value = "'707"
$ws.Range("D17").Value = "'6573657.56"
$ws.Range("C21").Value = "'236"
$ws.Range("D21").Value = "'744139.14"
$ws.Range("C23").Value = "'469"
$ws.Range("D23").Value = "'3586044.75"
$ws.Range("C24").Value = "'211"
$ws.Range("D24").Value = "'1396479.21"
$ws.Range("C34").Value = "'882"
$ws.Range("D34").Value = "'6849402.60"
$ws.Range("C52").Value = "'798"
$ws.Range("D52").Value = "'5186204.95"
$ws.Range("C59").Value = "'6819"
$ws.Range("D59").Value = "'34930560.31"
$ws.Range("C69").Value = "'250"
$ws.Range("D69").Value = "'741482.59"
$ws.Range("C80").Value = "'455"
$ws.Range("D80").Value = "'1479752.96"
$ws.Range("C82").Value = "'1278"
$ws.Range("D82").Value = "'10205367.06"
$ws.Range("C83").Value = "'671"
$ws.Range("D83").Value = "'4611277.80"
$ws.Range("C94").Value = "'270"
$ws.Range("D94").Value = "'738250.00"
$ws.Range("C96").Value = "'656"
$ws.Range("D96").Value = "'4504105.68"
$ws.Range("C104").Value = "'1705"
$ws.Range("D104").Value = "'9801214.95"
$ws.Range("C106").Value = "'1666"
$ws.Range("D106").Value = "'9105315.02"
$ws.Range("C108").Value = "'84"
$ws.Range("D108").Value = "'407961.23"
